$wb = $excel.ActiveWorkbook

# Sheet 1: "Metadata"
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B2").Value = "http://hl7.org/fhir/us/fhir-directory-query/ValueSet/SpecialtiesVS"
$wsMeta.Range("B8").Value = "2021-12-17T13:53:37-05:00"

# Sheet 2: "Include ValueSets"
$wsInclude1 = $wb.Worksheets.Item("Include ValueSets")
$wsInclude1.Range("A2").Value = "http://hl7.org/fhir/us/fhir-directory-query/ValueSet/IndividualAndGroupSpecialtiesVS"

# Sheet 3: "Include ValueSets 2"
$wsInclude2 = $wb.Worksheets.Item("Include ValueSets 2")
$wsInclude2.Range("A2").Value = "http://hl7.org/fhir/us/fhir-directory-query/ValueSet/NonIndividualSpecialtiesVS"
